$d = $word.ActiveDocument
$wordNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# -----------------------------------------------------------------
# Edit 1: move the "_GoBack" bookmark so it sits right after the
# first "dia anterior" paragraph (it currently only exists at the
# very end of the document).
# -----------------------------------------------------------------

# Remove the bookmark from its old location (end of document, after
# "...en el tipo invitado como en el superadmin.")
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-add it right after the run that reads "dia anterior" (first
# occurrence). A collapsed range placed exactly at the end of a
# paragraph's text cannot be targeted directly, so a one-character
# sentinel is appended, the bookmark is inserted in front of it
# (a perfectly safe, non-boundary position), and the sentinel is
# then deleted again -- leaving the bookmark exactly where the
# sentinel used to start, i.e. immediately after the real text.
$findRng = $d.Content
$findRng.Find.Execute("día anterior", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $findRng.End
$sentinelRng = $d.Range($endPos, $endPos)
$sentinelRng.InsertAfter("X")
$bmRng = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRng)
$sentinelRng2 = $d.Range($endPos, $endPos + 1)
$sentinelRng2.Delete()

# -----------------------------------------------------------------
# Edit 2: "Se puede medir: " -> split across several runs (simulating
# Word's spell/grammar-check markup around "medirhoy") and the
# paragraph alignment changes from justified to centered.
# -----------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("Se puede medir: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Alignment = 1

$xml = "<w:p $wordNs>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'>Se puede </w:t></w:r>" + `
  "<w:proofErr w:type='spellStart'/>" + `
  "<w:proofErr w:type='gramStart'/>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>medir</w:t></w:r>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>hoy</w:t></w:r>" + `
  "<w:proofErr w:type='spellEnd'/>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>" + `
  "<w:r><w:rPr><w:rFonts w:ascii='Arial' w:hAnsi='Arial' w:cs='Arial'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>:</w:t></w:r>" + `
  "<w:proofErr w:type='gramEnd'/>" + `
  "</w:p>"
$rng.InsertXML($xml)

# Remove the original "Se puede medir: " text, now left in place just
# before the newly-inserted runs.
$oldTextRng = $d.Content
$oldTextRng.Find.Execute("Se puede medir: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$oldTextRng.Text = ""
